$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.12586833333333
$ws.Range("H2").Value = 30.377605
$ws.Range("I2").Value = 0.9311967029481902
$ws.Range("J2").Value = 0.9311967029481902
$ws.Range("M2").Value = 31.618405
$ws.Range("N2").Value = 94.855215
$ws.Range("O2").Value = 0.8578613706944929
$ws.Range("P2").Value = 0.8578613706944929
$ws.Range("Q2").Value = 320.1638059400084
$ws.Range("R2").Value = 2881.474253460075
$ws.Range("S2").Value = 0.798837679977327
$ws.Range("T2").Value = 0.798837679977327

$ws.Range("G3").Value = 10.12586833333333
$ws.Range("H3").Value = 30.377605
$ws.Range("I3").Value = 0.9311967029481902
$ws.Range("J3").Value = 0.9311967029481902
$ws.Range("O3").Value = 0.08747555172986397
$ws.Range("P3").Value = 0.08747555172986396
$ws.Range("Q3").Value = 32.64688972515722
$ws.Range("R3").Value = 293.822007526415
$ws.Range("S3").Value = 0.08145694535942319
$ws.Range("T3").Value = 0.08145694535942317

$ws.Range("G4").Value = 10.12586833333333
$ws.Range("H4").Value = 30.377605
$ws.Range("I4").Value = 0.9311967029481902
$ws.Range("J4").Value = 0.9311967029481902
$ws.Range("M4").Value = 2.014730333333334
$ws.Range("N4").Value = 6.044191000000001
$ws.Range("O4").Value = 0.05466307757564324
$ws.Range("P4").Value = 0.05466307757564324
$ws.Range("Q4").Value = 20.40089408250612
$ws.Range("R4").Value = 183.608046742555
$ws.Range("S4").Value = 0.05090207761144014
$ws.Range("T4").Value = 0.05090207761144014

$ws.Range("I5").Value = 0.009287810103293732
$ws.Range("J5").Value = 0.009287810103293733
$ws.Range("M5").Value = 31.618405
$ws.Range("N5").Value = 94.855215
$ws.Range("O5").Value = 0.8578613706944929
$ws.Range("P5").Value = 0.8578613706944929
$ws.Range("Q5").Value = 3.19333243138
$ws.Range("R5").Value = 28.73999188242
$ws.Range("S5").Value = 0.00796765350596172
$ws.Range("T5").Value = 0.007967653505961722

$ws.Range("I6").Value = 0.009287810103293732
$ws.Range("J6").Value = 0.009287810103293733
$ws.Range("O6").Value = 0.08747555172986397
$ws.Range("P6").Value = 0.08747555172986396
$ws.Range("Q6").Value = 0.3256219779026666
$ws.Range("R6").Value = 2.930597801123999
$ws.Range("S6").Value = 0.0008124563131478241
$ws.Range("T6").Value = 0.0008124563131478241

$ws.Range("I7").Value = 0.009287810103293732
$ws.Range("J7").Value = 0.009287810103293733
$ws.Range("M7").Value = 2.014730333333334
$ws.Range("N7").Value = 6.044191000000001
$ws.Range("O7").Value = 0.05466307757564324
$ws.Range("P7").Value = 0.05466307757564324
$ws.Range("S7").Value = 0.0005077002841841883
$ws.Range("T7").Value = 0.0005077002841841884

$ws.Range("I8").Value = 0.05951548694851595
$ws.Range("J8").Value = 0.05951548694851596
$ws.Range("M8").Value = 31.618405
$ws.Range("N8").Value = 94.855215
$ws.Range("O8").Value = 0.8578613706944929
$ws.Range("P8").Value = 0.8578613706944929
$ws.Range("Q8").Value = 20.46259909800166
$ws.Range("R8").Value = 184.163391882015
$ws.Range("S8").Value = 0.05105603721120409
$ws.Range("T8").Value = 0.0510560372112041

$ws.Range("I9").Value = 0.05951548694851595
$ws.Range("J9").Value = 0.05951548694851596
$ws.Range("O9").Value = 0.08747555172986397
$ws.Range("P9").Value = 0.08747555172986396
$ws.Range("S9").Value = 0.005206150057292951
$ws.Range("T9").Value = 0.005206150057292951

$ws.Range("I10").Value = 0.05951548694851595
$ws.Range("J10").Value = 0.05951548694851596
$ws.Range("M10").Value = 2.014730333333334
$ws.Range("N10").Value = 6.044191000000001
$ws.Range("O10").Value = 0.05466307757564324
$ws.Range("P10").Value = 0.05466307757564324
$ws.Range("S10").Value = 0.00325329968001891
$ws.Range("T10").Value = 0.00325329968001891

